$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the two data records currently sitting in rows 2 and 3
# (row 2 becomes what row 3 held, and vice-versa) for the columns that
# are record-specific: Id/Taxonsorteringsordning/TaxonId/Artnamn/
# Vetenskapligt namn/Auktor/Antal/Ost/Nord. The rest of each row
# (Valideringsstatus, Lokalnamn, Noggrannhet, Lan, ... observatorer etc.)
# stays attached to its row number because both records share the same
# values there.

$numericCols = @("A","B","E","Q","R")
$textCols    = @("F","G","H")

foreach ($col in $numericCols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

foreach ($col in $textCols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

# Column I ("Antal") holds a text value "30" for the record now moving
# into row 2, and becomes blank for the record moving into row 3.
# Force it to stay text (it is stored as text in the source data) with
# a leading apostrophe so it isn't reinterpreted as a number.
$ws.Range("I2").Value = "'30"
$ws.Range("I3").Value = $null

# Column AF ("Bestamningsmetod") gains an empty placeholder cell on row
# 2 (mirroring the one the row-3 record used to carry) and loses it on
# row 3.
$ws.Range("AF2").Value = "'"
$ws.Range("AF3").Value = $null
